# Update "江西-漫展信息" workbook: refresh "want to go" counts on existing
# entries, and insert a newly-scraped event
# ("南昌·第三届龙年动漫展——庆端午贺高考专场") into the "展览" (Exhibition) and
# "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

function Update-Sheet($ws, $InsertRow, $FBefore, $NewRow, $FAfterShift) {

    # 1) Refresh "want to go" counts for rows that are NOT affected by the
    #    upcoming insert (they keep their row position).
    foreach ($r in $FBefore.Keys) {
        $ws.Range("F$r").Value = $FBefore[$r]
    }

    # 2) Insert a new blank row, pushing everything from $InsertRow down by one.
    $ws.Rows.Item($InsertRow).Insert()

    # Restore the bordered/bold index-column style on the new row's A cell
    # (mirrors the style used by every other row in column A).
    $srcIndexRow = $InsertRow - 1
    $ws.Range("A$srcIndexRow").Copy()
    $ws.Range("A$InsertRow").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0

    # 3) Populate the new row.
    $ws.Range("A$InsertRow").Value = ($InsertRow - 1)
    $ws.Range("B$InsertRow").NumberFormat = "@"
    $ws.Range("B$InsertRow").Value = $NewRow.B
    $ws.Range("C$InsertRow").Value = $NewRow.C
    $ws.Range("D$InsertRow").Value = $NewRow.D
    $ws.Range("E$InsertRow").NumberFormat = "@"
    $ws.Range("E$InsertRow").Value = $NewRow.E
    $ws.Range("F$InsertRow").Value = $NewRow.F
    $ws.Range("G$InsertRow").Value = $NewRow.G
    $ws.Range("H$InsertRow").Value = $NewRow.H
    $ws.Range("I$InsertRow").Value = $NewRow.I

    # 4) Fix up the running index in column A for every row that was pushed
    #    down by the insert: Insert() moves cell content but keeps the old
    #    values, so the "#" index column (A = row - 1) needs to be
    #    renumbered for rows InsertRow+1 .. LastRow.
    $lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1
    for ($r = $InsertRow + 1; $r -le $lastRow; $r++) {
        $ws.Range("A$r").Value = ($r - 1)
    }

    # 5) Refresh F/G on the rows that shifted down by one during the insert.
    foreach ($r in $FAfterShift.Keys) {
        $vals = $FAfterShift[$r]
        $ws.Range("F$r").Value = $vals[0]
        $ws.Range("G$r").Value = $vals[1]
    }
}

$newRowData = @{
    B = "2024-06-09"
    C = "南昌·第三届龙年动漫展——庆端午贺高考专场"
    D = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    E = "2024.06.09 10:00-06.10 18:00"
    F = 4
    G = 45
    H = "https://show.bilibili.com/platform/detail.html?id=85297"
    I = "//i1.hdslb.com/bfs/openplatform/202405/zBSAcG1V1714936299746.jpeg"
}

# ---------------------------------------------------------------------
# Sheet 1: "展览" (Exhibition) -- dimension A1:I20 -> A1:I21
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$fBefore1 = @{
    3  = 159
    5  = 9
    6  = 527
    7  = 1585
    8  = 7
    9  = 11
    10 = 1302
    11 = 117
}

$fAfterShift1 = @{
    14 = @(170, 50)
    16 = @(9, 55)
    17 = @(7, 55)
    18 = @(236, 64)
    20 = @(198, 40)
}

Update-Sheet $ws1 12 $fBefore1 $newRowData $fAfterShift1

# ---------------------------------------------------------------------
# Sheet 4: "全部类型" (All types) -- dimension A1:I21 -> A1:I22
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$fBefore4 = @{
    3  = 159
    5  = 9
    6  = 527
    7  = 1585
    9  = 7
    10 = 11
    11 = 1302
    12 = 117
}

$fAfterShift4 = @{
    15 = @(170, 50)
    17 = @(9, 55)
    18 = @(7, 55)
    19 = @(236, 64)
    21 = @(198, 40)
}

Update-Sheet $ws4 13 $fBefore4 $newRowData $fAfterShift4
